{"js": "// Insert three new bullet paragraphs describing Siege Analytics data engineering\n// achievements right after the \"Data Engineering and Infrastructure Architecture\"\n// sub-heading paragraph (and before the existing \"\u2022 Architect enterprise-scale...\"\n// bullet) in the PARTNER - Siege Analytics role.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the exact paragraph that introduces the Siege Analytics bullet list.\nlet anchorParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === \"Data Engineering and Infrastructure Architecture\") {\n    anchorParagraph = p;\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\n    \"Could not find paragraph 'Data Engineering and Infrastructure Architecture'\"\n  );\n}\n\nconst newBullets = [\n  \"\u2022 Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections\",\n  \"\u2022 Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government\",\n  \"\u2022 Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations\",\n];\n\n// Insert the three bullets immediately after the anchor paragraph, in order.\nlet insertAfter = anchorParagraph;\nfor (const bulletText of newBullets) {\n  insertAfter = insertAfter.insertParagraph(bulletText, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs describing Siege Analytics data engineering\n# achievements right after the \"Data Engineering and Infrastructure Architecture\"\n# sub-heading paragraph (and before the existing \"\u2022 Architect enterprise-scale...\"\n# bullet) in the PARTNER - Siege Analytics role.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Data Engineering and Infrastructure Architecture*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph 'Data Engineering and Infrastructure Architecture'\"\n}\n\n$newBullets = @(\n    \"\u2022 Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections\",\n    \"\u2022 Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government\",\n    \"\u2022 Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations\"\n)\n\n$cur = $target\nforeach ($bulletText in $newBullets) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $cur.Next()\n    $cur.Range.Text = $bulletText\n}\n"}
